$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab
$ws.Name = "BOM"

# Fix the bi-color LED naming in the BOM
$ws.Range("C14").Value = "LED 3mm Dome Bicolor"

# Select the edited cell as the active selection (matches resulting file state)
$ws.Range("C14").Select()
